$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to text
# (matching the original inline-string/text representation) by temporarily
# switching the cell to Text format, assigning the value, then restoring the
# default "Normal" style so no stray style index is left on the cell.

$ws.Range("D2").Value = "26.613.76"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.596.15"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").Value = "1.820.38"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").Value = "1.598.91"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "26.593.94"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("E21").Value = "  +4.84%  "

$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("E23").Value = "  -0.73%  "

$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.42"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("D34").Value = "1.282.61"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.621"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.38%  "

$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("E39").Value = "  +0.51%  "

$ws.Range("E40").Value = "  +19.47%  "

$ws.Range("E41").Value = "  +2.28%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D45").Value = "1.732.75"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("E48").Value = "  +4.31%  "

$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.10%  "
